$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric,
# so Excel does not auto-convert them to floats (matches source which stores them as text).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '24.658.33'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.657.38'
$ws.Range("E3").Value = '  -2.88%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = '320.65'
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").Value = '0.9989'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.3639'
$ws.Range("E7").Value = '  -2.74%  '
$ws.Range("D8").Value = '46.86'
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("D9").Value = '0.3262'
$ws.Range("E9").Value = '  -5.27%  '
$ws.Range("D10").Value = '1.133'
$ws.Range("E10").Value = '  -7.54%  '
$ws.Range("D11").Value = '0.07056'
$ws.Range("E11").Value = '  -6.43%  '
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '5.987'
$ws.Range("E13").Value = '  -5.60%  '
$ws.Range("D14").Value = '19.54'
$ws.Range("E14").Value = '  -7.86%  '
$ws.Range("D15").Value = '1.658.14'
$ws.Range("E15").Value = '  -2.98%  '
$ws.Range("D16").Value = '6.620'
$ws.Range("E16").Value = '  -6.30%  '
$ws.Range("D17").Value = '0.00001045'
$ws.Range("E17").Value = '  -7.64%  '
$ws.Range("D18").Value = '0.06628'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").Value = '0.9991'
$ws.Range("D20").Value = '78.89'
$ws.Range("E20").Value = '  -6.22%  '
$ws.Range("D21").Value = '5.923'
$ws.Range("E21").Value = '  -7.25%  '
$ws.Range("D22").Value = '15.75'
$ws.Range("E22").Value = '  -9.30%  '
$ws.Range("D23").Value = '12.64'
$ws.Range("E23").Value = '  -4.56%  '
$ws.Range("D24").Value = '24.658.68'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").Value = '2.466'
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").Value = '2.396'
$ws.Range("D27").Value = '147.97'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").Value = '18.60'
$ws.Range("E28").Value = '  -8.87%  '
$ws.Range("D29").Value = '1.843.58'
$ws.Range("E29").Value = '  -2.67%  '
$ws.Range("D30").Value = '1.211'
$ws.Range("E30").Value = '  -3.46%  '
$ws.Range("D31").Value = '125.46'
$ws.Range("E31").Value = '  -5.64%  '
$ws.Range("D32").Value = '4.073'
$ws.Range("E32").Value = '  -3.66%  '
$ws.Range("D33").Value = '5.839'
$ws.Range("E33").Value = '  -14.34%  '
$ws.Range("D34").Value = '0.08466'
$ws.Range("E34").Value = '  -3.83%  '
$ws.Range("D35").Value = '1.680'
$ws.Range("E35").Value = '  -5.63%  '
$ws.Range("D36").Value = '12.32'
$ws.Range("E36").Value = '  -10.82%  '
$ws.Range("D37").Value = '1.281'
$ws.Range("E37").Value = '  +0.72%  '
$ws.Range("D38").Value = '5.212'
$ws.Range("E38").Value = '  -7.41%  '
$ws.Range("D39").Value = '0.02236'
$ws.Range("E39").Value = '  -7.62%  '
$ws.Range("D40").Value = '0.06029'
$ws.Range("E40").Value = '  -9.58%  '
$ws.Range("D41").Value = '0.2076'
$ws.Range("E41").Value = '  -7.90%  '
$ws.Range("D42").Value = '8.211'
$ws.Range("E42").Value = '  -10.46%  '
$ws.Range("D43").Value = '0.9984'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = '0.5930'
$ws.Range("E44").Value = '  -8.78%  '
$ws.Range("D45").Value = '3.859'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("D46").Value = '12.83'
$ws.Range("E46").Value = '  -7.72%  '
$ws.Range("D47").Value = '0.5623'
$ws.Range("E47").Value = '  -8.93%  '
$ws.Range("D48").Value = '124.34'
$ws.Range("E48").Value = '  -4.01%  '
$ws.Range("D49").Value = '1.954'
$ws.Range("E49").Value = '  -7.92%  '
$ws.Range("D50").Value = '0.06971'
$ws.Range("E50").Value = '  -4.95%  '
$ws.Range("E51").Value = '  -3.11%  '
